$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 22.75000000000012
$ws.Range("H2").Value = 0.1608938545744234
$ws.Range("I2").Value = 0.1608938545744234
$ws.Range("L2").Value = 7.397655154538755
$ws.Range("M2").Value = "[-2.83839447934772, 17.63370478842523]"
$ws.Range("N2").Value = 0.1524433584489739
$ws.Range("O2").Value = 0.1524433584489739
$ws.Range("P2").Value = -0.9182633181663862
$ws.Range("Q2").Value = "[-4.006395436109506, 2.169868799776734]"
$ws.Range("R2").Value = 0.5522438347941148
$ws.Range("S2").Value = 0.5522438347941148
$ws.Range("T2").Value = 13.59083071794687
$ws.Range("U2").Value = "[8.200003829740513, 18.981657606153227]"
$ws.Range("V2").Value = [double]"7.09748791516418e-06"
$ws.Range("W2").Value = [double]"7.09748791516418e-06"
$ws.Range("X2").Value = 3.324824824824844
$ws.Range("Y2").Value = -7.856606606606645
$ws.Range("Z2").Value = 14.50625625625633

$ws.Range("F3").Value = 22.75000000000012
$ws.Range("H3").Value = 0.3570822783852973
$ws.Range("I3").Value = 0.3570822783852973
$ws.Range("L3").Value = 6.169301499949657
$ws.Range("M3").Value = "[-5.648090892454154, 17.98669389235347]"
$ws.Range("N3").Value = 0.2986574807592417
$ws.Range("O3").Value = 0.2986574807592417
$ws.Range("P3").Value = -0.8427896207828471
$ws.Range("Q3").Value = "[-3.9812375369816593, 2.295658295415965]"
$ws.Range("R3").Value = 0.5912705217958192
$ws.Range("S3").Value = 0.5912705217958192
$ws.Range("T3").Value = 17.9295470251978
$ws.Range("U3").Value = "[11.857577137450548, 24.00151691294505]"
$ws.Range("V3").Value = [double]"3.738171523437472e-07"
$ws.Range("W3").Value = [double]"3.738171523437472e-07"
$ws.Range("X3").Value = 3.051551551551565
$ws.Range("Y3").Value = -8.312062062062108
$ws.Range("Z3").Value = 14.41516516516524

$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 22.75000000000012
$ws.Range("H4").Value = 0.005872650503561228
$ws.Range("I4").Value = 0.005872650503561228
$ws.Range("L4").Value = 9.682074346823898
$ws.Range("M4").Value = "[1.6879265754646973, 17.676222118183098]"
$ws.Range("N4").Value = 0.01871909385962267
$ws.Range("O4").Value = 0.01871909385962267
$ws.Range("P4").Value = -0.767315923399309
$ws.Range("Q4").Value = "[-1.4591581494150798, -0.07547369738353815]"
$ws.Range("R4").Value = 0.0305091774412769
$ws.Range("S4").Value = 0.0305091774412769
$ws.Range("T4").Value = 13.04389140939376
$ws.Range("U4").Value = "[8.961466208759806, 17.126316610027718]"
$ws.Range("V4").Value = [double]"7.035641025332495e-08"
$ws.Range("W4").Value = [double]"7.035641025332495e-08"
$ws.Range("X4").Value = 2.778278278278293
$ws.Range("Y4").Value = 0.273273273273273
$ws.Range("Z4").Value = 5.283283283283313

$ws.Range("F5").Value = 22.75000000000012
$ws.Range("H5").Value = 0.1242231856775644
$ws.Range("I5").Value = 0.1242231856775644
$ws.Range("L5").Value = 6.891305662339375
$ws.Range("M5").Value = "[-1.6638151432166648, 15.446426467895416]"
$ws.Range("N5").Value = 0.1117046160061141
$ws.Range("O5").Value = 0.1117046160061141
$ws.Range("P5").Value = -0.1761052938949232
$ws.Range("Q5").Value = "[-1.8805529598065025, 1.528342372016656]"
$ws.Range("R5").Value = 0.8360904816678305
$ws.Range("S5").Value = 0.8360904816678305
$ws.Range("T5").Value = 11.99754233546003
$ws.Range("U5").Value = "[7.185360374114387, 16.809724296805676]"
$ws.Range("V5").Value = [double]"8.563748999490528e-06"
$ws.Range("W5").Value = [double]"8.563748999490528e-06"
$ws.Range("X5").Value = 0.637637637637642
$ws.Range("Y5").Value = -5.533783783783811
$ws.Range("Z5").Value = 6.809059059059095

$ws.Range("F6").Value = 22.75000000000012
$ws.Range("H6").Value = 0.6955894788646446
$ws.Range("I6").Value = 0.6955894788646446
$ws.Range("L6").Value = 2.452425061734066
$ws.Range("M6").Value = "[-5.139150725416588, 10.04400084888472]"
$ws.Range("N6").Value = 0.5185827610086979
$ws.Range("O6").Value = 0.5185827610086979
$ws.Range("P6").Value = -0.2264210921506153
$ws.Range("Q6").Value = "[-3.358579533567467, 2.905737349266236]"
$ws.Range("R6").Value = 0.8848891126036755
$ws.Range("S6").Value = 0.8848891126036755
$ws.Range("T6").Value = 11.85298475344798
$ws.Range("U6").Value = "[7.670982851742025, 16.034986655153944]"
$ws.Range("V6").Value = [double]"8.439077423538777e-07"
$ws.Range("W6").Value = [double]"8.439077423538777e-07"
$ws.Range("X6").Value = 0.8198198198198199
$ws.Range("Y6").Value = -10.52102102102108
$ws.Range("Z6").Value = 12.16066066066072

$ws.Range("F7").Value = 22.75000000000012
$ws.Range("H7").Value = 0.3234218969790599
$ws.Range("I7").Value = 0.3234218969790599
$ws.Range("L7").Value = 3.723422295301761
$ws.Range("M7").Value = "[-3.1106092504148535, 10.557453841018376]"
$ws.Range("N7").Value = 0.2783237456943048
$ws.Range("O7").Value = 0.2783237456943048
$ws.Range("P7").Value = 0.4339737599553457
$ws.Range("Q7").Value = "[-2.6918952066795425, 3.559842726590234]"
$ws.Range("R7").Value = 0.7810472666126282
$ws.Range("S7").Value = 0.7810472666126282
$ws.Range("T7").Value = 10.35541525560169
$ws.Range("U7").Value = "[6.8273950092346425, 13.88343550196873]"
$ws.Range("V7").Value = [double]"4.22062769178666e-07"
$ws.Range("W7").Value = [double]"4.22062769178666e-07"
$ws.Range("X7").Value = 21.17867867867879
$ws.Range("Y7").Value = 9.86061061061067
$ws.Range("Z7").Value = 32.49674674674691

$ws.Range("F8").Value = 22.75000000000012
$ws.Range("H8").Value = 0.1011819617080455
$ws.Range("I8").Value = 0.1011819617080455
$ws.Range("L8").Value = 8.6056100818633
$ws.Range("M8").Value = "[-2.1868886015451423, 19.398108765271743]"
$ws.Range("N8").Value = 0.1152728433059789
$ws.Range("O8").Value = 0.1152728433059789
$ws.Range("P8").Value = 0.786184347745194
$ws.Range("Q8").Value = "[-2.251631971942234, 3.824000667432622]"
$ws.Range("R8").Value = 0.6047497332690928
$ws.Range("S8").Value = 0.6047497332690928
$ws.Range("T8").Value = 17.78022682283749
$ws.Range("U8").Value = "[12.225856279905937, 23.334597365769053]"
$ws.Range("V8").Value = [double]"6.750765169094564e-08"
$ws.Range("W8").Value = [double]"6.750765169094564e-08"
$ws.Range("X8").Value = 19.90340340340351
$ws.Range("Y8").Value = 8.904154154154199
$ws.Range("Z8").Value = 30.90265265265281

$ws.Range("B9").Value = 0
$ws.Range("F9").Value = 22.75000000000012
$ws.Range("H9").Value = 0.1288786036713645
$ws.Range("I9").Value = 0.1288786036713645
$ws.Range("L9").Value = 7.292285851420515
$ws.Range("M9").Value = "[-1.439537512713735, 16.024109215554766]"
$ws.Range("N9").Value = 0.09948580013951491
$ws.Range("O9").Value = 0.09948580013951491
$ws.Range("P9").Value = 1.289342330302117
$ws.Range("Q9").Value = "[-0.5975001042863477, 3.176184764890581]"
$ws.Range("R9").Value = 0.1755380131577329
$ws.Range("S9").Value = 0.1755380131577329
$ws.Range("T9").Value = 15.34904369016907
$ws.Range("U9").Value = "[10.445701516658513, 20.252385863679624]"
$ws.Range("V9").Value = [double]"1.100252504659949e-07"
$ws.Range("W9").Value = [double]"1.100252504659949e-07"
$ws.Range("X9").Value = 18.08158158158168
$ws.Range("Y9").Value = 11.24974974974981
$ws.Range("Z9").Value = 24.91341341341354

$ws.Range("F10").Value = 24.18000000000034
$ws.Range("H10").Value = 0.4506184291159443
$ws.Range("I10").Value = 0.4506184291159443
$ws.Range("L10").Value = 3.932265290353818
$ws.Range("M10").Value = "[-3.910915007655714, 11.77544558836335]"
$ws.Range("N10").Value = 0.3179935178290099
$ws.Range("O10").Value = 0.3179935178290099
$ws.Range("P10").Value = 1.465447624197041
$ws.Range("Q10").Value = "[-1.6038160694001942, 4.534711317794276]"
$ws.Range("R10").Value = 0.3413636485688027
$ws.Range("S10").Value = 0.3413636485688027
$ws.Range("T10").Value = 11.60994900380652
$ws.Range("U10").Value = "[7.202741134401471, 16.01715687321157]"
$ws.Range("V10").Value = [double]"3.30257510205989e-06"
$ws.Range("W10").Value = [double]"3.30257510205989e-06"
$ws.Range("X10").Value = 18.54042042042068
$ws.Range("Y10").Value = 6.728768768768862
$ws.Range("Z10").Value = 30.3520720720725

$ws.Range("F11").Value = 24.18000000000034
$ws.Range("H11").Value = 0.1122300730878619
$ws.Range("I11").Value = 0.1122300730878619
$ws.Range("L11").Value = 6.663122517625089
$ws.Range("M11").Value = "[-1.2694103064040014, 14.59565534165418]"
$ws.Range("N11").Value = 0.09759987970549333
$ws.Range("O11").Value = 0.09759987970549333
$ws.Range("P11").Value = 1.314500229429963
$ws.Range("Q11").Value = "[-0.4465527095192714, 3.075553168379197]"
$ws.Range("R11").Value = 0.1397262156781771
$ws.Range("S11").Value = 0.1397262156781771
$ws.Range("T11").Value = 14.2876326411805
$ws.Range("U11").Value = "[9.850966816869944, 18.72429846549106]"
$ws.Range("V11").Value = [double]"5.911970513139408e-08"
$ws.Range("W11").Value = [double]"5.911970513139408e-08"
$ws.Range("X11").Value = 19.1213213213216
$ws.Range("Y11").Value = 12.34414414414432
$ws.Range("Z11").Value = 25.89849849849887

$ws.Range("F12").Value = 24.18000000000034
$ws.Range("H12").Value = 0.181566059668636
$ws.Range("I12").Value = 0.181566059668636
$ws.Range("L12").Value = 6.272607741280742
$ws.Range("M12").Value = "[-3.0093770211179063, 15.554592503679391]"
$ws.Range("N12").Value = 0.1802604879167324
$ws.Range("O12").Value = 0.1802604879167324
$ws.Range("P12").Value = 0.8993948938205012
$ws.Range("Q12").Value = "[-2.220184598032427, 4.018974385673429]"
$ws.Range("R12").Value = 0.5643537271829175
$ws.Range("S12").Value = 0.5643537271829175
$ws.Range("T12").Value = 10.53232886825884
$ws.Range("U12").Value = "[5.789118213784512, 15.275539522733176]"
$ws.Range("V12").Value = [double]"5.20596432345144e-05"
$ws.Range("W12").Value = [double]"5.20596432345144e-05"
$ws.Range("X12").Value = 20.71879879879909
$ws.Range("Y12").Value = 8.713513513513638
$ws.Range("Z12").Value = 32.72408408408455

$ws.Range("F13").Value = 24.18000000000034
$ws.Range("H13").Value = 0.126356726970106
$ws.Range("I13").Value = 0.126356726970106
$ws.Range("L13").Value = 6.810850249433962
$ws.Range("M13").Value = "[-2.237348402642972, 15.859048901510896]"
$ws.Range("N13").Value = 0.1364939529278044
$ws.Range("O13").Value = 0.1364939529278044
$ws.Range("P13").Value = 0.77360539818127
$ws.Range("Q13").Value = "[-2.364842518017542, 3.912053314380082]"
$ws.Range("R13").Value = 0.6219848030444899
$ws.Range("S13").Value = 0.6219848030444899
$ws.Range("T13").Value = 11.85229075143732
$ws.Range("U13").Value = "[7.196380824075124, 16.508200678799515]"
$ws.Range("V13").Value = [double]"6.016141459896929e-06"
$ws.Range("W13").Value = [double]"6.016141459896929e-06"
$ws.Range("X13").Value = 21.20288288288318
$ws.Range("Y13").Value = 9.124984984985117
$ws.Range("Z13").Value = 33.28078078078125

$ws.Range("B14").Value = 1
$ws.Range("F14").Value = 24.18000000000034
$ws.Range("H14").Value = 0.0009478808550952289
$ws.Range("I14").Value = 0.0009478808550952289
$ws.Range("L14").Value = 11.91560573613224
$ws.Range("M14").Value = "[3.9914529553871, 19.83975851687738]"
$ws.Range("N14").Value = 0.004059000666485124
$ws.Range("O14").Value = 0.004059000666485124
$ws.Range("P14").Value = 1.10065808684327
$ws.Range("Q14").Value = "[0.42139481039142357, 1.7799213632951174]"
$ws.Range("R14").Value = 0.002104617734991265
$ws.Range("S14").Value = 0.002104617734991265
$ws.Range("T14").Value = 12.90118085649254
$ws.Range("U14").Value = "[8.694056512334633, 17.108305200650456]"
$ws.Range("V14").Value = [double]"1.70854116632313e-07"
$ws.Range("W14").Value = [double]"1.70854116632313e-07"
$ws.Range("X14").Value = 19.94426426426455
$ws.Range("Y14").Value = 17.33021021021046
$ws.Range("Z14").Value = 22.55831831831864

$ws.Range("F15").Value = 24.18000000000034
$ws.Range("H15").Value = 0.6342776803213539
$ws.Range("I15").Value = 0.6342776803213539
$ws.Range("L15").Value = 3.47271700381166
$ws.Range("M15").Value = "[-6.0139887322458545, 12.959422739869174]"
$ws.Range("N15").Value = 0.4647758147136285
$ws.Range("O15").Value = 0.4647758147136285
$ws.Range("P15").Value = -0.1509473947670772
$ws.Range("Q15").Value = "[-3.2893953109658898, 2.9875005214317354]"
$ws.Range("R15").Value = 0.9232592340297312
$ws.Range("S15").Value = 0.9232592340297312
$ws.Range("T15").Value = 14.3795235762128
$ws.Range("U15").Value = "[9.24109465863058, 19.51795249379503]"
$ws.Range("V15").Value = [double]"1.07886105382704e-06"
$ws.Range("W15").Value = [double]"1.07886105382704e-06"
$ws.Range("X15").Value = 0.580900900900911
$ws.Range("Y15").Value = -11.49699699699716
$ws.Range("Z15").Value = 12.65879879879898

